$wb = $excel.ActiveWorkbook

# --- NewAcc1: append two new balance rows, update dimension/selection ---
$wsNewAcc1 = $wb.Worksheets.Item("NewAcc1")

# Copy the formatting of the existing data row (A3:B3) down into the two
# new rows before writing values, so the new cells pick up style s="10".
$wsNewAcc1.Range("A3:B3").Copy()
[void]$wsNewAcc1.Range("A4:B5").PasteSpecial(-4122) # xlPasteFormats

$wsNewAcc1.Range("A4").Value = 60182
$wsNewAcc1.Range("B4").Value = 2000
$wsNewAcc1.Range("A5").Value = 38278
$wsNewAcc1.Range("B5").Value = 3000

# --- CustomerDetails: move the active selection ---
$wsCustomerDetails = $wb.Worksheets.Item("CustomerDetails")
[void]$wsCustomerDetails.Range("I2").Select()

# --- Add the new BalEnq sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "BalEnq"

# Reuse NewAcc1's header/data formatting (s="12" header, s="10" data) for
# the new sheet's cells. (NewAcc1's own selection is already being moved
# elsewhere in this script, so using it as a copy source adds no extra
# unwanted diffs there.)
$wsNewAcc1.Range("A1").Copy()
[void]$newSheet.Range("A1").PasteSpecial(-4122)
$wsNewAcc1.Range("A2").Copy()
[void]$newSheet.Range("A2:A3").PasteSpecial(-4122)

$newSheet.Range("A1").Value = "Accno"
$newSheet.Range("A2").Value = 89666
$newSheet.Range("A3").Value = 34189
[void]$newSheet.Range("D3").Select()

# --- Restore NewAcc1 as the active/selected sheet+cell ---
[void]$wsNewAcc1.Activate()
[void]$wsNewAcc1.Range("E6").Select()
